$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.153.88'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '2.242.52'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  +0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '246.54'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.82%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.629'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.58%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '76.04'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +4.84%  '
$ws.Range("E8").Value = '  -0.07%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.630'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.23%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '40.16'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.68%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0949'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.95%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.22'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.78%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.104'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").Value = '2.578.87'
$ws.Range("E14").Value = '  -2.15%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '14.86'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.96%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.860'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.98%  '
$ws.Range("D17").Value = '2.246.68'
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("D18").Value = '42.134.34'
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("D19").Value = '0.0₃0977'
$ws.Range("E19").Value = '  -2.20%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.16'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.35%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '71.49'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.61%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -4.18%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '231.48'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("E24").Value = '  -0.12%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '3.73'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -4.86%  '
$ws.Range("E26").Value = '  -4.34%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.30%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.08'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +9.92%  '
$ws.Range("E29").Value = '  -1.36%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '168.24'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.44%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '20.51'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.81%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.0848'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +5.55%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '32.52'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.95%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -7.03%  '
$ws.Range("E35").Value = '  -0.08%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -5.33%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.79'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.93%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0296'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.64%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '13.28'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -5.50%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.49%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -6.14%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '117.58'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +21.24%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.202'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -4.98%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '60.05'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.23%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.70'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -6.22%  '
$ws.Range("E46").Value = '  -2.73%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -4.63%  '
$ws.Range("E49").Value = '  -1.63%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '4.27'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -12.92%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.60%  '
